$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "2025/12/03 02:13"
$ws.Range("B9").Value = "-"
$ws.Range("C9").Value = "-"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "-"
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "-"
